$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - "Shape 76" (hashtag text box): merge the " #" run and the
# "Clemson #" run into a single " #Clemson #" run (same formatting, so
# re-writing the combined span collapses the two <a:r> elements into one).
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(2)
$tr1 = $shp1.TextFrame.TextRange
$tr1.Characters(17, 11).Text = " #Clemson #"

# ---------------------------------------------------------------------------
# Slides 2-6 - the title ("Shape 82") and subtitle ("Shape 83") textboxes
# move down by the same amount (167748 EMU = 13.2085 pt).
# ---------------------------------------------------------------------------
for ($i = 2; $i -le 6; $i++) {
    $sl = $p.Slides.Item($i)

    $title = $sl.Shapes.Item(1)
    $title.Top = 384283 / 12700

    $subtitle = $sl.Shapes.Item(2)
    $subtitle.Top = 972969 / 12700
}

# ---------------------------------------------------------------------------
# Slide 2 - "Shape 85" body text box: merge several adjacent same-format
# runs (a side effect of retyping across the run boundary).
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(3)
$tr2 = $shp2.TextFrame.TextRange

$tr2.Characters(2, 52).Text = "Spoken at various conferences including SaltConf15, "
$tr2.Characters(63, 34).Text = " Summit 2015 (Tokyo), SaltConf16, "
$tr2.Characters(105, 31).Text = " conference 2016 (Germany) etc."
$tr2.Characters(290, 173).Text = "Few of the significant contributions made include making Salt compatible with Python 3, creating salt-cloud VMware driver, creating DNS/ASAM/Spacewalk runners, creating ZFS/"

# ---------------------------------------------------------------------------
# Slide 5 - "Shape 83" subtitle text box: merge "what is it? " and
# "why was it created?" runs.
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(2)
$tr5 = $shp5.TextFrame.TextRange
$tr5.Characters(1, 31).Text = "what is it? why was it created?"

# ---------------------------------------------------------------------------
# Slide 7 - "Shape 76" (hashtag text box): same merge as slide 1.
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(4)
$tr7 = $shp7.TextFrame.TextRange
$tr7.Characters(17, 11).Text = " #Clemson #"
